$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary_counts")

# Insert a new row at position 13, pushing existing rows 13-14 down to 14-15
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new data
$ws.Range("A13").Value = "Number of events with both any university response coding and any police coding"
$ws.Range("B13").Value = 63
